# Applies the 06-11-2023 02:45 scraper update to the Costa Rica Primera
# Division 2023-2024 sheet:
#   1) Rows 91 and 92 had their match data (columns F:V) swapped - the
#      "Sporting San Jose vs Guanacasteca" match and the "Cartagines vs
#      Saprissa" match traded places (the A:E identifying columns stay put).
#   2) A new match row (108) was appended: Grecia 0 x 0 Herediano.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the F:V content of rows 91 and 92 -----------------------------

$row91 = @($ws.Range("F91").Value2, $ws.Range("G91").Value2, $ws.Range("H91").Value2, `
    $ws.Range("I91").Value2, $ws.Range("J91").Value2, $ws.Range("K91").Value2, `
    $ws.Range("L91").Value2, $ws.Range("M91").Value2, $ws.Range("N91").Value2, `
    $ws.Range("O91").Value2, $ws.Range("P91").Value2, $ws.Range("Q91").Value2, `
    $ws.Range("R91").Value2, $ws.Range("S91").Value2, $ws.Range("T91").Value2, `
    $ws.Range("U91").Value2, $ws.Range("V91").Value2)

$row92 = @($ws.Range("F92").Value2, $ws.Range("G92").Value2, $ws.Range("H92").Value2, `
    $ws.Range("I92").Value2, $ws.Range("J92").Value2, $ws.Range("K92").Value2, `
    $ws.Range("L92").Value2, $ws.Range("M92").Value2, $ws.Range("N92").Value2, `
    $ws.Range("O92").Value2, $ws.Range("P92").Value2, $ws.Range("Q92").Value2, `
    $ws.Range("R92").Value2, $ws.Range("S92").Value2, $ws.Range("T92").Value2, `
    $ws.Range("U92").Value2, $ws.Range("V92").Value2)

$cols = @("F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

for ($idx = 0; $idx -lt $cols.Length; $idx++) {
    $ws.Range($cols[$idx] + "91").Value = $row92[$idx]
    $ws.Range($cols[$idx] + "92").Value = $row91[$idx]
}

# --- 2) Append the new match as row 108 -------------------------------------

# Copy formatting (styles) from the last existing data row (107) down to the
# new row so the "Indice" cell keeps the bold/bordered style and the date
# cell keeps its date number format.
$ws.Range("A107:V107").Copy()
$ws.Range("A108:V108").PasteSpecial(-4122)

$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "costa-rica"
$ws.Range("C108").Value = "primera-division"
$ws.Range("D108").Value = "2023-2024"
$ws.Range("E108").Value = 45235.91666666666
$ws.Range("F108").Value = "Grecia"
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = "Herediano"
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 4.76
$ws.Range("K108").Value = "31/10/2023 14:42"
$ws.Range("L108").Value = 5.23
$ws.Range("M108").Value = "05/11/2023 21:51"
$ws.Range("N108").Value = 3.92
$ws.Range("O108").Value = "31/10/2023 14:42"
$ws.Range("P108").Value = 4.23
$ws.Range("Q108").Value = "05/11/2023 21:51"
$ws.Range("R108").Value = 1.68
$ws.Range("S108").Value = "31/10/2023 14:42"
$ws.Range("T108").Value = 1.61
$ws.Range("U108").Value = "05/11/2023 21:51"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/grecia-herediano/8zXhflBU/"

"applied"
